# Generate Report for Handback
# Updates the localization-status workbook: marks handback as complete,
# stamps handback datetimes, and fills in the "Latest Target File" /
# "Latest Handback File" columns (F, G) on the per-language sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: update status cells ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("F2").Value = "4257e57a-322f-4761-94e0-60c09395b319.md"
$wsZh.Range("G2").Value = "4257e57a-322f-4761-94e0-60c09395b319.d8fecea4e7fa830f2985f63fa6f491c97599d269.zh-cn.xlf"
$wsZh.Range("F3").Value = "4257e57a-322f-4761-94e0-60c09395b319.md"
$wsZh.Range("G3").Value = "4257e57a-322f-4761-94e0-60c09395b319.d8fecea4e7fa830f2985f63fa6f491c97599d269.zh-cn.xlf"

$wsZh.Range("H2").Value = "2016-03-20 04:52:24"
$wsZh.Range("H3").Value = "2016-03-20 04:52:24"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("F2").Value = "4257e57a-322f-4761-94e0-60c09395b319.md"
$wsDe.Range("G2").Value = "4257e57a-322f-4761-94e0-60c09395b319.d8fecea4e7fa830f2985f63fa6f491c97599d269.de-de.xlf"
$wsDe.Range("F3").Value = "4257e57a-322f-4761-94e0-60c09395b319.md"
$wsDe.Range("G3").Value = "4257e57a-322f-4761-94e0-60c09395b319.d8fecea4e7fa830f2985f63fa6f491c97599d269.de-de.xlf"

$wsDe.Range("H2").Value = "2016-03-20 04:52:29"
$wsDe.Range("H3").Value = "2016-03-20 04:52:29"
